$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") value corrections per repulled data
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 5
$ws.Range("F6").Value = -11
$ws.Range("F12").Value = 2
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = -5
